$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the scratch/staging data that lived in columns F:H (rows 1-26).
# That data held draft variable names (F), descriptions (G) and categories (H)
# which are being promoted into the main A:C dimension table below.
$ws.Range("F1:H26").Clear()

# Append the promoted dimension rows to the bottom of the main table in
# columns A (variable), B (description) and C (category), starting at row 177.

$ws.Range("A177").Value = 'SEXADJB'
$ws.Range("B177").Value = 'This field represents the number of levels applied under Part B of the Repeat and Dangerous Sex Offender adjustment (§4B1.5). Use with MONSEXB.'
$ws.Range("C177").Value = 'calc'
$ws.Rows(177).RowHeight = 28.8

$ws.Range("A178").Value = 'STA1_1'
$ws.Range("B178").Value = 'Title, Section, and Subsection number, respectively of the first statutes for each count of conviction. For example, STA1_123 is the 1st statute for the 123rd count of conviction. USSC records up to 3 statutes for each count. These statute fields are not ordered; therefore, the first
statute is not necessarily the most serious count. Use with the variable NOCOUNT. See NWSTAT1-
NWSTATX to find UNIQUE statutes from all fields.'
$ws.Range("C178").Value = 'crime'
$ws.Rows(178).RowHeight = 100.8

$ws.Range("A179").Value = 'STA2_1'
$ws.Range("B179").Value = 'Title, Section, and Subsection number, respectively of the second statutes for each count of conviction. For example, STA2_123 is the 2nd statute for the 123rd count of conviction. USSC records up to 3 statutes for each count. These statute fields are not ordered; therefore, the first
statute is not necessarily the most serious count. Use with the variable NOCOUNT. See NWSTAT1-
NWSTATX to find UNIQUE statutes from all fields.'
$ws.Range("C179").Value = 'crime'
$ws.Rows(179).RowHeight = 100.8

$ws.Range("A180").Value = 'DAFROM1'
$ws.Range("B180").Value = 'The low end of the range of drug weight for first through nth drug types (DRUGTYP1-DRUGTYPX). This field is used when an exact amount was not specified in DRGAM1-DRGAMX but a range was provided. The weights are reported in several units of measure (UNIT1-UNITX).'
$ws.Range("C180").Value = 'drug_related'
$ws.Rows(180).RowHeight = 57.6

$ws.Range("A181").Value = 'DAFROM2'
$ws.Range("B181").Value = 'The low end of the range of drug weight for first through nth drug types (DRUGTYP1-DRUGTYPX). This field is used when an exact amount was not specified in DRGAM1-DRGAMX but a range was provided. The weights are reported in several units of measure (UNIT1-UNITX).'
$ws.Range("C181").Value = 'drug_related'
$ws.Rows(181).RowHeight = 57.6

$ws.Range("A182").Value = 'DRGAM1'
$ws.Range("B182").Value = 'Drug amount of the first through nth drug types (DRUGTYP1-DRUGTYPX) for which the defendant was held responsible. Often when weights are missing it is because both parties agreed to a Base Offense Level and the documents do not specify a corresponding drug amount. The weights are reported in several units of measure (UNIT1-UNITX), see WGT1- WGTX for the gram weight equivalency. For drug range amounts, see DAFROM1-DAFROMX and DATO1-DATOX.'
$ws.Range("C182").Value = 'drug_related'
$ws.Rows(182).RowHeight = 86.4

$ws.Range("A183").Value = 'DRGAM2'
$ws.Range("B183").Value = 'Drug amount of the first through nth drug types (DRUGTYP1-DRUGTYPX) for which the defendant was held responsible. Often when weights are missing it is because both parties agreed to a Base Offense Level and the documents do not specify a corresponding drug amount. The weights are reported in several units of measure (UNIT1-UNITX), see WGT1- WGTX for the gram weight equivalency. For drug range amounts, see DAFROM1-DAFROMX and DATO1-DATOX.'
$ws.Range("C183").Value = 'drug_related'
$ws.Rows(183).RowHeight = 86.4

$ws.Range("A184").Value = 'DATO1'
$ws.Range("B184").Value = 'The high end of the range of drug weight for first through nth drug types (DRUGTYP1-DRUGTYPX) when an exact amount was not specified in DRGAM1-DRGAMX. The weights are reported in several units of measure (UNIT1-UNITX).'
$ws.Range("C184").Value = 'drug_related'
$ws.Rows(184).RowHeight = 43.2

$ws.Range("A185").Value = 'DATO2'
$ws.Range("B185").Value = 'The high end of the range of drug weight for first through nth drug types (DRUGTYP1-DRUGTYPX) when an exact amount was not specified in DRGAM1-DRGAMX. The weights are reported in several units of measure (UNIT1-UNITX).'
$ws.Range("C185").Value = 'drug_related'
$ws.Rows(185).RowHeight = 43.2

$ws.Range("A186").Value = 'IMMIMIN'
$ws.Range("B186").Value = 'Represents mandatory minimum sentence (in months) associated with 8§1324 (Immigration).'
$ws.Range("C186").Value = 'mand_mins'
$ws.Rows(186).RowHeight = 28.8

$ws.Range("A187").Value = 'LOSS1'
$ws.Range("B187").Value = 'The dollar amount of loss for which the offender is held responsible. Amounts  are rounded off to the nearest whole dollar. Loss is often used in various economic crime guidelines to determine either the base offense level or levels of an SOC. Cases not involving dollar loss (ex. Drug cases) are coded as zero for the amount. Use variable NOCOMP to determine how many guideline computations are present in each case. All guideline variables available'
$ws.Range("C187").Value = 'crime'
$ws.Rows(187).RowHeight = 86.4

$ws.Range("A188").Value = 'LOSS2'
$ws.Range("B188").Value = 'The dollar amount of loss for which the offender is held responsible. Amounts  are rounded off to the nearest whole dollar. Loss is often used in various economic crime guidelines to determine either the base offense level or levels of an SOC. Cases not involving dollar loss (ex. Drug cases) are coded as zero for the amount. Use variable NOCOMP to determine how many guideline computations are present in each case. All guideline variables available'
$ws.Range("C188").Value = 'crime'
$ws.Rows(188).RowHeight = 86.4

$ws.Range("A189").Value = 'MARRIED'
$ws.Range("B189").Value = 'Marital status of offender. This field is available FY1999-FY2003.'
$ws.Range("C189").Value = 'bio'

$ws.Range("A190").Value = 'MNTHDEPT'
$ws.Range("B190").Value = 'The difference in months between the guideline minimum (GLMIN) and the sentence length, including alternatives and probation as 0 months or incarceration (SENSPCAP). Only present for above and below range cases where the GLMIN is greater than 0 and less than life and the sentence is not life. Large values of GLMIN/MAX above sentencing table excluded from calcs. Field available FY2018-present.'
$ws.Range("C190").Value = 'stat'
$ws.Rows(190).RowHeight = 72

$ws.Range("A191").Value = 'MWEIGHT'
$ws.Range("B191").Value = 'The marijuana weight equivalency, in grams, of all the drug types coded. This variable is not missing if ANY of the individual marijuana equivalency weights are available (i.e., if the case involves two drugs and one type has weight available and one type does not have the weight available, then MWEIGHT will not be missing).'
$ws.Range("C191").Value = 'drug_related'
$ws.Rows(191).RowHeight = 57.6

$ws.Range("A192").Value = 'MWGT1'
$ws.Range("B192").Value = 'The marijuana weight equivalency, in grams, of the first through nth drug types (DRUGTYP1-DRUGTYPX)'
$ws.Range("C192").Value = 'drug_related'
$ws.Rows(192).RowHeight = 28.8

$ws.Range("A193").Value = 'MWGT2'
$ws.Range("B193").Value = 'The marijuana weight equivalency, in grams, of the first through nth drug types (DRUGTYP1-DRUGTYPX)'
$ws.Range("C193").Value = 'drug_related'
$ws.Rows(193).RowHeight = 28.8

$ws.Range("A194").Value = 'NWSTAT1'
$ws.Range("B194").Value = 'Title, Section, and Subsection number of the UNIQUE statutes for each case generated from all of the statute fields (STAT1_1 thru STAT3_XX). Ex. If a case has a total of 5 counts involving 4 counts of 21:841 and one count of 18:924C then NWSTAT1 will be 21:841, NWSTAT2 will be 18:924C and NWSTAT3-X will be
missing/inapplicable. These statute fields are not ordered; therefore, the first statute is not necessarily the most serious count. Use with NOUSTAT.'
$ws.Range("C194").Value = 'crime'
$ws.Rows(194).RowHeight = 86.4

$ws.Range("A195").Value = 'NWSTAT2'
$ws.Range("B195").Value = 'Title, Section, and Subsection number of the UNIQUE statutes for each case generated from all of the statute fields (STAT1_1 thru STAT3_XX). Ex. If a case has a total of 5 counts involving 4 counts of 21:841 and one count of 18:924C then NWSTAT1 will be 21:841, NWSTAT2 will be 18:924C and NWSTAT3-X will be
missing/inapplicable. These statute fields are not ordered; therefore, the first statute is not necessarily the most serious count. Use with NOUSTAT.'
$ws.Range("C195").Value = 'crime'
$ws.Rows(195).RowHeight = 86.4

$ws.Range("A196").Value = 'OFFTYPSB'
$ws.Range("B196").Value = 'Primary offense type variable used in the FY2010 through FY2017 Sourcebook tables (this field replaces OFFTYPE2 although OFFTYPE2 is still available on the datafile through FY2017). This variable is based on the count of conviction with the highest statutory maximum (in case of a tie, the count with the highest statutory minimum is used). Note that since the primary offense type is derived from statutes of conviction it may not match up logically with the primary guideline (or any of the guidelines applied). See OFFTYPE2 for offense types used in USSC Sourcebook FY1999-FY2009. See OFFGUIDE for offense types used in USSC Sourcebook FY2018-present. Field available FY2010-FY2017.'
$ws.Range("C196").Value = 'crime'
$ws.Rows(196).RowHeight = 129.6

$ws.Range("A197").Value = 'PCNTDEPT'
$ws.Range("B197").Value = 'The percent difference between the guideline minimum (GLMIN) and the sentence length, including alternatives and probation as 0 months or incarceration (SENSPCAP). Only present for above and below range cases where the GLMIN is greater than 0 and less than life and the sentence is not life. Large values of GLMIN/MAX above sentencing table excluded from calcs. Field available FY2018-present.'
$ws.Range("C197").Value = 'stat'
$ws.Rows(197).RowHeight = 72

$ws.Range("A198").Value = 'BOOKERCD'
$ws.Range("B198").Value = 'Assigns cases to one of the 12 post-Booker reporting categories based on relationship between the sentence and guideline range and the reason(s) given
for being outside of the range. '
$ws.Range("C198").Value = 'stat'
$ws.Rows(198).RowHeight = 43.2

$ws.Range("A199").Value = 'SENTTOT0'
$ws.Range("B199").Value = 'The total prison sentence (excluding months of alternative confinement), in months, with zeros (probation). Missing cases are set to "ꞏ". This field includes sentences of time imposed, time served, and §5G1.3 credit. See Appendix B in this codebook for more information about USSC sentencing variables. Field available FY1999-FY2017. See also SENTTCAP.'
$ws.Range("C199").Value = 'outcome'
$ws.Rows(199).RowHeight = 72

# Re-fit column C (categories) now that longer values have been added to it.
$ws.Columns("C").AutoFit()

# Leave the selection where the authoring session ended up (just past the
# last new row).
[void]$ws.Range("C200").Select()
